$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q6").Value = 88.157250792756912
Write-Output $ws.Range("Q6").Value
